# Auto-generated COM-interop script applying the "Add household surplus computations" commit.
$wb = $excel.ActiveWorkbook

# ===== Sheet: Summary =====
$ws = $wb.Worksheets.Item('Summary')
# Copy row 6 formatting down into the new rows 7:10
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range('B6').Value = 2896154.550748563
$ws.Range('A7').Value = 'Wasted Prosumer Surplus'
$ws.Range('B7').Value = 342857.5917216506
$ws.Range('A8').Value = 'Total Wasted Prosumer Surplus'
$ws.Range('B8').Value = 342857.5917216506
$ws.Range('A9').Value = 'Unmet Demand'
$ws.Range('B9').Value = 29364.14959724166
$ws.Range('A10').Value = 'Total Unmet Demand'
$ws.Range('B10').Value = 29364.14959724166

# ===== Sheet: Costs and Revenues =====
$ws = $wb.Worksheets.Item('Costs and Revenues')
$ws.Range('B4').Value = 516729.2344596348
$ws.Range('B6').Value = 400618.8526664367
$ws.Range('C4').Value = 514932.0834136077
$ws.Range('C6').Value = 406416.0037124637
$ws.Range('D4').Value = 513132.4944200165
$ws.Range('D6').Value = 408215.5927060549
$ws.Range('E2').Value = 938737.9016593837
$ws.Range('E4').Value = 501136.9470735459
$ws.Range('E6').Value = 101234.3895858377
$ws.Range('F2').Value = 959349.9008549106
$ws.Range('F3').Value = 20800
$ws.Range('F4').Value = 511935.5369542995
$ws.Range('F5').Value = 34552.359
$ws.Range('F6').Value = 392062.0049006111
$ws.Range('G2').Value = 959349.9008549106
$ws.Range('G3').Value = 4000
$ws.Range('G4').Value = 510118.6542394872
$ws.Range('G6').Value = 410678.8876154234
$ws.Range('H2').Value = 959349.9008549105
$ws.Range('H4').Value = 508299.2497951454
$ws.Range('H6').Value = 416498.2920597651
$ws.Range('I2').Value = 959349.9008549105
$ws.Range('I4').Value = 506477.3054326685
$ws.Range('I6').Value = 418320.236422242
$ws.Range('J4').Value = 489991.0232228274
$ws.Range('J6').Value = 132868.5327648077
$ws.Range('K2').Value = 961571.7997246224
$ws.Range('K3').Value = 48000
$ws.Range('K4').Value = 504151.9865617157
$ws.Range('K6').Value = 374531.1781629067
$ws.Range('L2').Value = 969896.833024751
$ws.Range('L3').Value = 18400
$ws.Range('L4').Value = 507265.0393204461
$ws.Range('L6').Value = 407829.916704305
$ws.Range('M2').Value = 969896.833024751
$ws.Range('M4').Value = 505406.1812562778
$ws.Range('M6').Value = 428088.7747684733
$ws.Range('N2').Value = 969896.8330247512
$ws.Range('N4').Value = 503544.650463314
$ws.Range('N6').Value = 429950.3055614372
$ws.Range('O2').Value = 916492.3349876349
$ws.Range('O4').Value = 470341.4757097192
$ws.Range('O6').Value = 192517.8812779157
$ws.Range('P2').Value = 838926.15105909
$ws.Range('P4').Value = 423288.3152576661
$ws.Range('P5').Value = 25388.838
$ws.Range('P6').Value = 390248.9978014239

# ===== Sheet: Installed Capacities =====
$ws = $wb.Worksheets.Item('Installed Capacities')
$ws.Range('F2').Value = 411
$ws.Range('P2').Value = 302

# ===== Sheet: Added Capacities =====
$ws = $wb.Worksheets.Item('Added Capacities')
$ws.Range('F2').Value = 26
$ws.Range('G2').Value = 5
$ws.Range('K2').Value = 60
$ws.Range('L2').Value = 23

# ===== Sheet: Retired Capacities =====
$ws = $wb.Worksheets.Item('Retired Capacities')
$ws.Range('K2').Value = 26
$ws.Range('L2').Value = 5
$ws.Range('P2').Value = 60

# ===== Sheet: DG Dispatch =====
$ws = $wb.Worksheets.Item('DG Dispatch')
$ws.Range('B14').Value = 411
$ws.Range('B44').Value = 302
$ws.Range('B45').Value = 302
$ws.Range('C14').Value = 411
$ws.Range('C44').Value = 302
$ws.Range('C45').Value = 302
$ws.Range('D14').Value = 410.3391557398498
$ws.Range('D44').Value = 302
$ws.Range('D45').Value = 302
$ws.Range('E44').Value = 302
$ws.Range('E45').Value = 302
$ws.Range('F44').Value = 302
$ws.Range('F45').Value = 302
$ws.Range('G14').Value = 410.8327491714383
$ws.Range('G44').Value = 302
$ws.Range('G45').Value = 302
$ws.Range('H14').Value = 411
$ws.Range('H44').Value = 302
$ws.Range('H45').Value = 302
$ws.Range('Q16').Value = 411
$ws.Range('Q46').Value = 302
$ws.Range('R16').Value = 411
$ws.Range('R45').Value = 302
$ws.Range('R46').Value = 302
$ws.Range('S45').Value = 302
$ws.Range('S46').Value = 302
$ws.Range('T14').Value = 411
$ws.Range('T44').Value = 302
$ws.Range('T45').Value = 302
$ws.Range('U14').Value = 411
$ws.Range('U44').Value = 302
$ws.Range('U45').Value = 302
$ws.Range('V14').Value = 411
$ws.Range('V15').Value = 411
$ws.Range('V44').Value = 302
$ws.Range('V45').Value = 302
$ws.Range('W14').Value = 411
$ws.Range('W15').Value = 411
$ws.Range('W44').Value = 302
$ws.Range('W45').Value = 302
$ws.Range('X14').Value = 411
$ws.Range('X15').Value = 411
$ws.Range('X44').Value = 302
$ws.Range('X45').Value = 302
$ws.Range('Y14').Value = 411
$ws.Range('Y44').Value = 302
$ws.Range('Y45').Value = 302

# ===== Sheet: Unmet Demand =====
$ws = $wb.Worksheets.Item('Unmet Demand')
$ws.Range('B14').Value = 70.99931295557451
$ws.Range('B44').Value = 179.9993129555745
$ws.Range('B45').Value = 82.55655664632661
$ws.Range('C14').Value = 38.47457824299391
$ws.Range('C44').Value = 147.4745782429939
$ws.Range('C45').Value = 59.09991244551929
$ws.Range('D14').Value = 0
$ws.Range('D44').Value = 108.3391557398498
$ws.Range('D45').Value = 45.93768689770263
$ws.Range('E44').Value = 102.3632896068686
$ws.Range('E45').Value = 40.67209722191262
$ws.Range('F44').Value = 102.8896287080119
$ws.Range('F45').Value = 37.63624233787687
$ws.Range('G14').Value = 0
$ws.Range('G44').Value = 108.8327491714383
$ws.Range('G45').Value = 27.52519625238585
$ws.Range('H14').Value = 24.02773927029563
$ws.Range('H44').Value = 133.0277392702956
$ws.Range('H45').Value = 43.22842014979517
$ws.Range('Q16').Value = 111.1821235684552
$ws.Range('Q46').Value = 220.1821235684552
$ws.Range('R16').Value = 310.1956210454637
$ws.Range('R45').Value = 51.67054165050009
$ws.Range('R46').Value = 419.1956210454637
$ws.Range('S45').Value = 38.140588939824
$ws.Range('S46').Value = 107.541226054864
$ws.Range('T14').Value = 118.6191915811053
$ws.Range('T44').Value = 227.6191915811053
$ws.Range('T45').Value = 86.53401876295709
$ws.Range('U14').Value = 233.5217529288726
$ws.Range('U44').Value = 342.5217529288726
$ws.Range('U45').Value = 97.68869740971195
$ws.Range('V14').Value = 218.8510241668239
$ws.Range('V15').Value = 3.510667191520156
$ws.Range('V44').Value = 327.8510241668239
$ws.Range('V45').Value = 112.5106671915202
$ws.Range('W14').Value = 227.3734759809475
$ws.Range('W15').Value = 21.37314290982852
$ws.Range('W44').Value = 336.3734759809475
$ws.Range('W45').Value = 130.3731429098285
$ws.Range('X14').Value = 181.2818334606677
$ws.Range('X15').Value = 8.862739445387547
$ws.Range('X44').Value = 290.2818334606677
$ws.Range('X45').Value = 117.8627394453875
$ws.Range('Y14').Value = 100.3174326828064
$ws.Range('Y44').Value = 209.3174326828064
$ws.Range('Y45').Value = 97.39139276613435

Write-Output "edit complete"